# Update "想去人数" (want-to-go count, column F) figures on the "展览" and
# "全部类型" sheets to the freshly scraped totals.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 1127
$ws1.Range("F8").Value  = 11500
$ws1.Range("F9").Value  = 4342
$ws1.Range("F10").Value = 28
$ws1.Range("F12").Value = 19
$ws1.Range("F14").Value = 2532
$ws1.Range("F15").Value = 1082
$ws1.Range("F16").Value = 128
$ws1.Range("F17").Value = 30
$ws1.Range("F18").Value = 2207
$ws1.Range("F19").Value = 180
$ws1.Range("F20").Value = 509
$ws1.Range("F21").Value = 11292
$ws1.Range("F22").Value = 11199

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value  = 1127
$ws4.Range("F8").Value  = 11500
$ws4.Range("F9").Value  = 4342
$ws4.Range("F10").Value = 28
$ws4.Range("F12").Value = 19
$ws4.Range("F14").Value = 2532
$ws4.Range("F16").Value = 1082
$ws4.Range("F17").Value = 128
$ws4.Range("F18").Value = 30
$ws4.Range("F19").Value = 2207
$ws4.Range("F20").Value = 180
$ws4.Range("F21").Value = 509
$ws4.Range("F22").Value = 11292
$ws4.Range("F23").Value = 11199
